# Apply the "Add files via upload" revision to the 牛顿环 (Newton's rings)
# workbook:
#   1. Update the GitHub repo link in column A (row 21) from
#      "fuck-university-physics-experiments" to "fuck-nku-physics-experiments".
#      (Written first so the shared-string table gets the same append order
#      the original author's save produced.)
#   2. Reword the red-cell instructional note in A3 to the longer wording.
#   3. Add two helper cells (J17/K17) next to the wavelength-uncertainty
#      result in row 17: a text-formatted echo of I17, and that value + 1.
#   4. Clear the trailing Touhou quote cell (row 23), leaving it blank but
#      still styled.
#   5. Move the active selection to J20, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the posted-on link text.
$ws.Range("A21").Value = "Posted on https://github.com/Axolyz/fuck-nku-physics-experiments."

# 2. Reword the red-cell instructions.
$ws.Range("A3").Value = "红色格子：填入你的实验数据，如本身自带数据请更改"

# 3. New helper cells beside I17.
$ws.Range("J17").Formula = '=TEXT(I17,"0.0000")'
$ws.Range("K17").Formula = "=J17+1"

# 4. Clear the closing quote, keep the row/style in place.
$ws.Range("A23").Value = $null

# 5. Restore the saved selection.
$ws.Range("J20").Select()
